{"js": "// Applies the \"Desenvolvimento da Home e altera\u00e7\u00e3o dos requisitos da\n// documenta\u00e7\u00e3o\" edits to the Peneirando requirements document.\n//\n// Five textual spots change (two of them are pure run-splits around a\n// word Word's spell-checker flagged, i.e. the rendered text is\n// unchanged there; the other three actually change the wording):\n//\n//  1. \", como de Vin\u00edcius Jr. e Rodrygo, ambos transferidos...\" \u2013 text\n//     unchanged (Word just re-split the run around \"Rodrygo\").\n//  2. \"a, mesmo em um cen\u00e1rio pr\u00e9-pandemia\" \u2013 text unchanged (re-split\n//     around \"pr\u00e9\").\n//  3. \"...(local, data, hor\u00e1rio e requisitos);\" -> \"...(Clube, local,\n//     data, hor\u00e1rio e informa\u00e7\u00f5es gerais);\"\n//  4. \"Na aba \u201cSobre n\u00f3s\u201d, deve haver...\" -> \"Ainda na aba \u201cPeneiras\u201d\n//     deve haver...\"\n//  5. \"a aba \u201cmeu perfil\u201d fica habilitada...\" -> \"a aba \u201cPerfil\u201d fica\n//     habilitada...\"\n//\n// Since 1 & 2 do not change any visible character, they are no-ops for\n// the purposes of this script (Office.js has no public API to stamp\n// <w:proofErr/> spell-check markers, and doing so would not alter the\n// document's text content).\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 3. \"Peneiras\" bullet: add \"Clube, \" and replace \"requisitos\" with\n//    \"informa\u00e7\u00f5es gerais\".\nawait replaceOnce(\n  \"com os dados da peneira (local, data, hor\u00e1rio e requisitos);\",\n  \"com os dados da peneira (Clube, local, data, hor\u00e1rio e informa\u00e7\u00f5es gerais);\"\n);\n\n// 4. \"Sobre n\u00f3s\" bullet becomes a continuation of the \"Peneiras\" topic.\nawait replaceOnce(\n  \"Na aba \\u201cSobre n\\u00f3s\\u201d, deve haver um breve texto sobre nosso trabalho e miss\\u00e3o\",\n  \"Ainda na aba \\u201cPeneiras\\u201d deve haver um breve texto sobre nosso trabalho e miss\\u00e3o\"\n);\n\n// 5. \"meu perfil\" -> \"Perfil\".\nawait replaceOnce(\n  \"a aba \\u201cmeu perfil\\u201d fica habilitada\",\n  \"a aba \\u201cPerfil\\u201d fica habilitada\"\n);\n", "ps1": "# Applies the \"Desenvolvimento da Home e altera\u00e7\u00e3o dos requisitos da\n# documenta\u00e7\u00e3o\" edits to the Peneirando requirements document.\n#\n# Three spots in the requirements list actually change wording:\n#   3. \"...(local, data, hor\u00e1rio e requisitos);\"\n#      -> \"...(Clube, local, data, hor\u00e1rio e informa\u00e7\u00f5es gerais);\"\n#   4. \"Na aba \u201cSobre n\u00f3s\u201d, deve haver...\" -> \"Ainda na aba \u201cPeneiras\u201d\n#      deve haver...\"\n#   5. \"a aba \u201cmeu perfil\u201d fica habilitada...\" -> \"a aba \u201cPerfil\u201d fica\n#      habilitada...\"\n#\n# (Two other hunks in the upstream diff only re-split existing runs\n# around \"Rodrygo\" / \"pr\u00e9\" \u2014 Word stamping fresh <w:proofErr/> spell-\n# check markers while the author typed nearby \u2014 with no change to the\n# visible characters, so there is nothing to replay for them here.)\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once([string]$findText, [string]$replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute(\n        $findText,    # FindText\n        $true,        # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        0,            # Wrap (wdFindStop)\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        1             # Replace (wdReplaceOne)\n    )\n    if (-not $ok) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# 3. \"Peneiras\" bullet: add \"Clube, \" and replace \"requisitos\" with\n#    \"informa\u00e7\u00f5es gerais\".\nReplace-Once `\n    \"com os dados da peneira (local, data, hor\u00e1rio e requisitos);\" `\n    \"com os dados da peneira (Clube, local, data, hor\u00e1rio e informa\u00e7\u00f5es gerais);\"\n\n# 4. \"Sobre n\u00f3s\" bullet becomes a continuation of the \"Peneiras\" topic.\nReplace-Once `\n    \"Na aba \u201cSobre n\u00f3s\u201d, deve haver um breve texto sobre nosso trabalho e miss\u00e3o\" `\n    \"Ainda na aba \u201cPeneiras\u201d deve haver um breve texto sobre nosso trabalho e miss\u00e3o\"\n\n# 5. \"meu perfil\" -> \"Perfil\".\nReplace-Once `\n    \"a aba \u201cmeu perfil\u201d fica habilitada\" `\n    \"a aba \u201cPerfil\u201d fica habilitada\"\n"}
